$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1664724431887483
$ws.Range("C2").Value = 1.464482041830579
$ws.Range("D2").Value = 6.721781159715255
$ws.Range("E2").Value = 2.592639805240068
$ws.Range("F2").Value = 2.648175453457832
$ws.Range("G2").Value = 22
